# Generate Report for Handoff
# - Updates the status/timestamp for the "0ced5dab..." file row to reflect
#   it being ready for handoff again.
# - Removes the row tracking the "4bf7fb59..." file (handled/closed out),
#   on all three sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Update row 2 (0ced5dab...) with the new status/date.
$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"
$ov.Range("D2").Value = "2016-03-19 08:01:19"

# Remove row 3 (4bf7fb59...) entirely.
$ov.Rows.Item(3).Delete()

# Rebuild the hyperlinks collection (row delete leaves stale entries).
$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8354f1bca090a8a831688b991a19fe21b273aa98/e2e/0ced5dab-da16-40e8-91a1-cba26979ea92.md", [Type]::Missing, [Type]::Missing, "0ced5dab-da16-40e8-91a1-cba26979ea92.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Update row 2 (0ced5dab...) with the new status/date.
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("E2").Value = "2016-03-19 08:01:09"

# Remove row 3 (4bf7fb59...) entirely.
$zh.Rows.Item(3).Delete()

# Rebuild the hyperlinks collection (row delete leaves stale entries).
$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8354f1bca090a8a831688b991a19fe21b273aa98/e2e/0ced5dab-da16-40e8-91a1-cba26979ea92.md", [Type]::Missing, [Type]::Missing, "0ced5dab-da16-40e8-91a1-cba26979ea92.md")
$zh.Hyperlinks.Add($zh.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/45f2d397115b39f889319b401dce60a2a4761bc2/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/0ced5dab-da16-40e8-91a1-cba26979ea92.7f4df12ac097546f8afaa9e3310709186c06e91a.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "0ced5dab-da16-40e8-91a1-cba26979ea92.7f4df12ac097546f8afaa9e3310709186c06e91a.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-zhcn-fly/blob/b60d25d06b1ce713634a1ced46722ff131f74899/e2e/0ced5dab-da16-40e8-91a1-cba26979ea92.md", [Type]::Missing, [Type]::Missing, "0ced5dab-da16-40e8-91a1-cba26979ea92.md")
$zh.Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/3797b5228c2eb94c38303b2e27b0f697b51b5c98/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/0ced5dab-da16-40e8-91a1-cba26979ea92.7f4df12ac097546f8afaa9e3310709186c06e91a.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "0ced5dab-da16-40e8-91a1-cba26979ea92.7f4df12ac097546f8afaa9e3310709186c06e91a.zh-cn.xlf")

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Update row 2 (0ced5dab...) with the new status/date.
$de.Range("C2").Value = "Ready for handoff"
$de.Range("E2").Value = "2016-03-19 08:01:19"

# Remove row 3 (4bf7fb59...) entirely.
$de.Rows.Item(3).Delete()

# Rebuild the hyperlinks collection (row delete leaves stale entries).
$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8354f1bca090a8a831688b991a19fe21b273aa98/e2e/0ced5dab-da16-40e8-91a1-cba26979ea92.md", [Type]::Missing, [Type]::Missing, "0ced5dab-da16-40e8-91a1-cba26979ea92.md")
$de.Hyperlinks.Add($de.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/36c95a8bfe3dd7dde96c4580c823d639ac5c6494/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/0ced5dab-da16-40e8-91a1-cba26979ea92.7f4df12ac097546f8afaa9e3310709186c06e91a.de-de.xlf", [Type]::Missing, [Type]::Missing, "0ced5dab-da16-40e8-91a1-cba26979ea92.7f4df12ac097546f8afaa9e3310709186c06e91a.de-de.xlf")
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest-dede-fly/blob/132019f2d41cfbe58298004a88b2b3046c45f1c1/e2e/0ced5dab-da16-40e8-91a1-cba26979ea92.md", [Type]::Missing, [Type]::Missing, "0ced5dab-da16-40e8-91a1-cba26979ea92.md")
$de.Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a9f492919e2ae5f9d580e9013b83b81082a7c12a/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/0ced5dab-da16-40e8-91a1-cba26979ea92.7f4df12ac097546f8afaa9e3310709186c06e91a.de-de.xlf", [Type]::Missing, [Type]::Missing, "0ced5dab-da16-40e8-91a1-cba26979ea92.7f4df12ac097546f8afaa9e3310709186c06e91a.de-de.xlf")

Write-Host "Done"
